$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new timesheet entry on row 14 (previously a blank placeholder row
# that already had the shared formulas in columns D:G but no input data)
$ws.Range("A14").Value = 45577
$ws.Range("B14").Value = 0.62708333333333333
$ws.Range("C14").Value = 0.73333333333333328

# Extend the week-3 summary ranges (row 4) so they roll up the new row 14
# entry along with rows 10-13
$ws.Range("M4").Formula = "=SUM(D10:D14)"
$ws.Range("N4").Formula = "=SUM(G10:G14)"

# Update the active selection to reflect where the user was last working
[void]$ws.Range("N5").Select()
